$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Title heading (appears twice - main heading and bold paragraph near the end)
Replace-Text "Play Jungle Gold for Free - Review of Jungle Gold Slot Game" "Play Jungle Gold for Free - Review of Unique Slot Game"

# "What we like" bullet list
Replace-Text "Stunning graphics in a cartoon-style jungle theme" "Unique game mechanics with expanding reels"
Replace-Text "Engaging and entertaining soundtrack" "Stunning graphics and immersive jungle theme"
Replace-Text "Lucrative bonus features increase chances of winning big payouts" "Engaging soundtrack that keeps players entertained"
Replace-Text "Betting range from 0.10 to 100 coins per spin" "Lucrative bonus features for big payouts"

# "What we don't like" bullet list
Replace-Text "Only a single line of five boxes, with two more revealed through bonus features" "Limited number of paylines with only a single line of five boxes"
Replace-Text "Lowest value symbols are differently colored and shaped gems" "Higher betting range may not be suitable for all players"

# Meta description (italic paragraph at the end)
Replace-Text "Play Jungle Gold for free with our review. Discover the pros and cons of Jungle Gold slot game, with its stunning graphics and lucrative bonus features." "Discover the features of Jungle Gold and play for free. Find out if it's worth your time!"
